$d = $word.ActiveDocument

# These character styles have <w:rPr> children that are out of wml.xsd
# schema order (color before b/i). Re-assigning the Bold/Italic property
# via the object model forces the run properties to be re-serialized in
# the canonical CT_RPr sequence (rFonts, b, bCs, i, iCs, ..., color, ...),
# fixing the OOXMLValidator Sch_UnexpectedElementContentExpectingComplex
# warning without changing any visible formatting.

$boldOnlyStyles = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnlyStyles) {
    $style = $d.Styles($styleName)
    $style.Font.Bold = $true
}

$italicOnlyStyles = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnlyStyles) {
    $style = $d.Styles($styleName)
    $style.Font.Italic = $true
}

$boldItalicStyles = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldItalicStyles) {
    $style = $d.Styles($styleName)
    $style.Font.Bold = $true
    $style.Font.Italic = $true
}
